# Update the "Soma" summary sheet with the real Slither/smartbugs analysis
# results (counts per vulnerability category), sorted descending by count.
# Column A keeps the original stable index of each category, column B is
# the category label, column C is the total count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, "reentrancy", 22),
    @(1, "arithmetic", 20),
    @(9, "Other", 16),
    @(4, "unchecked_low_calls", 10),
    @(0, "access_control", 7),
    @(7, "time_manipulation", 6),
    @(5, "bad_randomness", 4),
    @(2, "denial_service", 0),
    @(6, "front_running", 0),
    @(8, "short_addresses", 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
